$d = $word.ActiveDocument

# Locate the sentence boundary "...na outro. Resumidamente" so we can anchor the
# insertion precisely between "outro" and the following period.
$anchor = $d.Content
$found = $anchor.Find.Execute("outro. Resumidamente", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Insertion point is right after "outro" (5 characters into the match).
    $insertAt = $d.Range($anchor.Start + 5, $anchor.Start + 5)
    $insertAt.InsertAfter(" – não pode ser comprovado")

    # Make "não pode ser comprovado" bold (leave the " – " plain).
    $boldTarget = $d.Content
    $foundBold = $boldTarget.Find.Execute("não pode ser comprovado", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundBold) {
        $boldTarget.Bold = 1
    }
}
